# 1st changes of mifos to finflux
#
# The "Repayment schedule" sheet gets a new blank column inserted before the
# existing "Late" column (so: # | Days | Date | Paid Date | heading |
# Principal Due | Balance of Loan | Interest | Fees | Penalties | Due | Paid |
# In Advance | <new blank column> | Late | heading | Outstanding), and this
# sheet becomes the active tab/selection of the workbook (moved away from
# "Sheet1").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make "Repayment schedule" the active sheet (this also updates
# workbook.xml's bookViews/workbookView@activeTab and moves
# sheetView@tabSelected away from the previously active "Sheet1" tab).
$ws.Activate()

# Insert a new blank column at position N (14), pushing the existing
# N/O/P columns (Late / heading / Outstanding) one to the right.
$ws.Columns("N:N").Insert() | Out-Null

# Match the width of the new column to its left neighbour (column M,
# "In Advance") so the inserted column keeps a sensible custom width.
$ws.Columns("N:N").ColumnWidth = 9.83

# Reflect the new selection left behind on the sheet.
$ws.Range("M7").Select() | Out-Null
